$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'We don''t have much money now.'
$ws.Range("C5").Value = 'You can''t enter this building with dangerous things such as knife.'
$ws.Range("C34").Value = 'I can''t stop reading this book. It''s so interesting.'
$ws.Range("C38").Value = 'I''m looking for a dress similar to yours.'
$ws.Range("C45").Value = 'No, I usually get up quite late.'
$ws.Range("C46").Value = 'How often do you go skiing in every winter?'
$ws.Range("C47").Value = 'Neither could I. We have to work hard tonight.'
$ws.Range("C48").Value = 'So is my sister. The flu is going around.'
$ws.Range("C49").Value = 'No, she never cooks.'
$ws.Range("C50").Value = 'Don''t rush me. Eating slowly is healthier than eating quickly.'
$ws.Range("C55").Value = 'My friend didn''t have the courage to talk to the famous actor, and I didn''t, either.'
$ws.Range("C58").Value = 'I don''t eat much meat but I eat a lot of vegetable.'
$ws.Range("C64").Value = 'I don''t have much knowledge of the subject.'
$ws.Range("C66").Value = 'I want a new bicycle, but my father won''t buy me one.'
$ws.Range("C68").Value = 'I didn''t know where to go or what to do.'
$ws.Range("C75").Value = 'I don''t doubt that Mary will be happy to hear the news.'
$ws.Range("B78").Value = '私はお店に買い物をしている間に財布を無くしました。'
$ws.Range("C79").Value = 'Didn''t you wash your hands before you ate lunch?'
$ws.Range("C80").Value = 'You can''t use this facility unless you live in this city.'
$ws.Range("C81").Value = 'Let''s go swimming if it isn''t rainy tomorrow.'
$ws.Range("C84").Value = 'The doll house was so delicate that we couldn''t touch it.'
$ws.Range("C86").Value = 'Don''t make a sound, or you will wake up that dog.'
$ws.Range("C90").Value = 'I won''t tell it to you if you don''t promise to keep it a secret.'
$ws.Range("C92").Value = 'The ceiling was so high that I couldn''t change the light bulbs.'
$ws.Range("C93").Value = 'Get on this train, or you won''t get there by noon.'
$ws.Range("C95").Value = 'You should watch tonight''s game even if you aren''t interested in soccer.'
$ws.Range("C96").Value = 'This curry is so hot that I can''t eat all of it.'
$ws.Range("C102").Value = 'Let''s meet at the station.'
$ws.Range("C106").Value = 'It''s a nice day. Let''s walk through the park.'
$ws.Range("C107").Value = 'I couldn''t see the stage well because there was a tall man sitting in front of me.'
$ws.Range("C114").Value = 'My mother''s watch was made in Switzerland.'
$ws.Range("C117").Value = 'I can''t believe that paper is made from plants.'
$ws.Range("C118").Value = 'Let''s race to the top of the mountain.'
$ws.Range("C121").Value = 'We cannot enter the park until nine o''clock.'
$ws.Range("C123").Value = 'It is raining hard. I don''t want to go out of the building.'
$ws.Range("C126").Value = 'Let''s hope nothing like this will happen again.'
$ws.Range("C135").Value = 'Didn''t you meet Jane at the station?'
$ws.Range("C136").Value = 'He wasn''t born in Japan, was he?'
$ws.Range("C137").Value = 'I don''t have much time. How about you?'
$ws.Range("C139").Value = 'Why don''t you ask her about it on the phone?'
$ws.Range("C140").Value = 'Maki and Dina are good friend aren''t they?'
$ws.Range("C151").Value = 'Why isn''t Jack playing on this team?'
$ws.Range("C152").Value = 'The meeting does start until 10 o''clock, doesn''t it?'
$ws.Range("C153").Value = 'Why don''t we invite Bob and Naomi to the barbecue?'
$ws.Range("C164").Value = 'The boy said to the soccer player, " I have always wanted to talk to you."'
$ws.Range("C166").Value = 'The tourist said to the conductor, "Please tell me the way to the hotel."'
$ws.Range("C167").Value = 'The old man always says to us "How are you?"'
$ws.Range("C172").Value = 'The villager asked him why he hadn''t came to work the previous week?'
$ws.Range("C173").Value = 'The students said to the woman "We are going to stay at this hotel tonight.'
$ws.Range("C174").Value = 'I didn''t know you were a friend of the famous writer.'
$ws.Range("C175").Value = 'I didn''t know you were a friend of the famous writer.'
